# Add a new worksheet "serlogin" (an API test dataset for the login service)
# after the existing "login" sheet, and wire up selection / hyperlinks to
# match the authored workbook.

$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item("login")

# Select an entire column on the "login" sheet (matches the recorded
# pre-switch selection state) before we move away from it.
$loginSheet.Columns.Item(2).Select() | Out-Null

# Duplicate "login" so the new sheet starts life with the same look & feel
# (fonts / row height / margins / phonetic settings) as the original sheet,
# then wipe its contents and repurpose it.
$loginSheet.Copy($null, $loginSheet) | Out-Null
$newSheet = $wb.Worksheets.Item("login (2)")
$newSheet.Name = "serlogin"
$newSheet.Cells.Clear()

$url = "http://192.168.31.70:5000/login"

# Header row
$newSheet.Range("A1").Value = "id"
$newSheet.Range("B1").Value = "url"
$newSheet.Range("C1").Value = "body"
$newSheet.Range("D1").Value = "desc"
$newSheet.Range("E1").Value = "method"
$newSheet.Range("F1").Value = "expected"

# Row 2 - successful login
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = $url
$newSheet.Range("C2").Value = '{"username": "user1", "password": "password1"}'
$newSheet.Range("D2").Value = "登陆成功"
$newSheet.Range("E2").Value = "post"
$newSheet.Range("F2").Value = '{"message": "登录成功"}'

# Row 3 - wrong password
$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = $url
$newSheet.Range("C3").Value = '{"username": "user1", "password": "password"}'
$newSheet.Range("D3").Value = "密码错误"
$newSheet.Range("E3").Value = "post"
$newSheet.Range("F3").Value = '{"message": "用户名或密码错误"}'

# Row 4 - missing password
$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = $url
$newSheet.Range("C4").Value = '{"username": "user1"}'
$newSheet.Range("D4").Value = "缺少用户名或密码"
$newSheet.Range("E4").Value = "post"
$newSheet.Range("F4").Value = '{"message": "缺少用户名或密码"}'

# Row 5 - missing username
$newSheet.Range("A5").Value = 4
$newSheet.Range("B5").Value = $url
$newSheet.Range("C5").Value = '{"password": "password"}'
$newSheet.Range("D5").Value = "缺少用户名或密码"
$newSheet.Range("E5").Value = "post"
$newSheet.Range("F5").Value = '{"message": "缺少用户名或密码"}'

# Hyperlink the url cells out to the live endpoint (Excel auto-applies the
# built-in "Hyperlink" cell style, s="1", to each of these).
$newSheet.Hyperlinks.Add($newSheet.Range("B2"), $url)
$newSheet.Hyperlinks.Add($newSheet.Range("B3"), $url)
$newSheet.Hyperlinks.Add($newSheet.Range("B4"), $url)
$newSheet.Hyperlinks.Add($newSheet.Range("B5"), $url)

# Column widths
$newSheet.Columns.Item(2).ColumnWidth = 34.5546875
$newSheet.Columns.Item(3).ColumnWidth = 44.109375
$newSheet.Columns.Item(4).ColumnWidth = 19.88671875
$newSheet.Columns.Item(5).ColumnWidth = 13.5546875
$newSheet.Columns.Item(6).ColumnWidth = 32.6640625

# Row 2 (first data row) is taller to fit the wrapped JSON body text.
$newSheet.Rows.Item(2).RowHeight = 27

# Selection left on the new sheet, then make it the active tab.
$newSheet.Range("F7").Select() | Out-Null
$newSheet.Activate() | Out-Null
